$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (avoid Excel auto-converting numeric-looking strings like "1.038" into
# actual numbers) by forcing a Text number format while we write values,
# then restoring the default "Normal" style afterwards.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.803.05'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '1.859.06'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '1.038'
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '323.58'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D7").Value = '0.4414'
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("D8").Value = '0.3821'
$ws.Range("E8").Value = '  +1.64%  '
$ws.Range("D10").Value = '0.8879'
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("D11").Value = '21.63'
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '1.848.36'
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").Value = '5.545'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '6.746'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '0.07215'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("E16").Value = '  +3.96%  '
$ws.Range("D17").Value = '1.041'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '0.000009116'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '1.034'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").Value = '15.59'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").Value = '27.823.84'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").Value = '5.302'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = '11.27'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '2.101.77'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = '2.066'
$ws.Range("E25").Value = '  +6.22%  '
$ws.Range("D26").Value = '159.25'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").Value = '2.010'
$ws.Range("E28").Value = '  +3.56%  '
$ws.Range("D29").Value = '5.384'
$ws.Range("D30").Value = '118.79'
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("D31").Value = '0.09117'
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").Value = '1.220'
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("D33").Value = '0.7759'
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.027'
$ws.Range("E34").Value = '  +4.89%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '4.616'
$ws.Range("E35").Value = '  +2.14%  '
$ws.Range("D36").Value = '1.036'
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").Value = '1.158'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '0.01986'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").Value = '0.05327'
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").Value = '2.869'
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("D41").Value = '0.5221'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").Value = '6.961'
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").Value = '8.825'
$ws.Range("E44").Value = '  +2.69%  '
$ws.Range("D45").Value = '111.17'
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").Value = '10.84'
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").Value = '1.036'
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").Value = '0.06594'
$ws.Range("E48").Value = '  +3.02%  '
$ws.Range("D49").Value = '1.722'
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").Value = '0.4736'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").Value = '1.891'
$ws.Range("E51").Value = '  +0.18%  '

$priceRange.Style = "Normal"
